$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.972.71'
$ws.Range("E2").Value = '  -2.29%  '

$ws.Range("D3").Value = '2.664.94'
$ws.Range("E3").Value = '  -1.28%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '''595.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("D6").Value = '''163.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.19%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '''0.545'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.33%  '

$ws.Range("D9").Value = '2.662.81'
$ws.Range("E9").Value = '  -1.36%  '

$ws.Range("D10").Value = '''0.140'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.40%  '

$ws.Range("D11").Value = '''0.157'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").Value = '''0.356'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").Value = '''5.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.13%  '

$ws.Range("D14").Value = '''27.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.24%  '

$ws.Range("D15").Value = '3.164.67'

$ws.Range("D16").Value = '''0.0000181'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.86%  '

$ws.Range("D17").Value = '67.008.07'
$ws.Range("E17").Value = '  -2.22%  '

$ws.Range("D18").Value = '2.660.33'
$ws.Range("E18").Value = '  -1.43%  '

$ws.Range("D19").Value = '''11.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.46%  '

$ws.Range("D20").Value = '''360.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").Value = '''7.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.14%  '

$ws.Range("D22").Value = '''4.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.66%  '

$ws.Range("D23").Value = '''4.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.49%  '

$ws.Range("D24").Value = '''2.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.86%  '

$ws.Range("D25").Value = '''71.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.09%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").Value = '''9.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("E28").Value = '  -1.85%  '

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''1.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '''0.0000101'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.59%  '

$ws.Range("D31").Value = '''550.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.93%  '

$ws.Range("D32").Value = '''7.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.03%  '

$ws.Range("D33").Value = '''1.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.65%  '

$ws.Range("D34").Value = '''1.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.23%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '''1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.128'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.11%  '

$ws.Range("D37").Value = '''1.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.80%  '

$ws.Range("D38").Value = '''19.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.29%  '

$ws.Range("D39").Value = '''156.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.77%  '

$ws.Range("D40").Value = '''0.371'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.17%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = '''5.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.92%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.96%  '

$ws.Range("D43").Value = '''17.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.80%  '

$ws.Range("D46").Value = '''40.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.63%  '

$ws.Range("D47").Value = '0.0₆0297'
$ws.Range("E47").Value = '  -6.24%  '

$ws.Range("D48").Value = '''0.583'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.46%  '

$ws.Range("D49").Value = '''152.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.29%  '

$ws.Range("D50").Value = '''3.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.96%  '

$ws.Range("D51").Value = '''1.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.26%  '

